# robot framework folder updated
#
# 1. login_data: 3rd row email changes from Rani13@gmail.com -> Rani16@gmail.com
# 2. productsort_data and searchproduct_data swap names (and their data) while
#    staying in the same tab position (tab position / sheetId stay put, only
#    the sheet name and the cell content swap)
# 3. registration_data: sample email/password become Neeraja43@gmail.com / Neeraja123@AP
# 4. a brand new sheet "filterorder_data" is appended at the end - a copy of the
#    original login_data layout (email/password/username) using the old
#    Rani13@gmail.com row

$wb = $excel.ActiveWorkbook

# --- 1. login_data tweak -------------------------------------------------
$login = $wb.Worksheets.Item("login_data")
$login.Cells.Item(3, 1).Value = "Rani16@gmail.com"

# --- 2. swap productsort_data <-> searchproduct_data --------------------
$sortSheet = $wb.Worksheets.Item("productsort_data")
$searchSheet = $wb.Worksheets.Item("searchproduct_data")

# rename in place (no .Move -> tab position / sheetId stay where they are)
$sortSheet.Name = "tmp_swap_name"
$searchSheet.Name = "productsort_data"
$sortSheet.Name = "searchproduct_data"

# $sortSheet is now named "searchproduct_data" and must hold the product list
# $searchSheet is now named "productsort_data" and must hold the sort list

# shrink the (now) searchproduct_data sheet from 7 rows down to 3
$sortSheet.Range("A4:A7").EntireRow.Delete()

$productData = @("product_name", "Selenium Ruby", "Robot")
for ($i = 0; $i -lt $productData.Length; $i++) {
    $sortSheet.Cells.Item($i + 1, 1).Value = $productData[$i]
}

# grow the (now) productsort_data sheet from 3 rows up to 7
$sortData = @("sort_by", "menu_order", "popularity", "rating", "date", "price", "price-desc")
for ($i = 0; $i -lt $sortData.Length; $i++) {
    $searchSheet.Cells.Item($i + 1, 1).Value = $sortData[$i]
}

# --- 3. registration_data sample values ----------------------------------
$reg = $wb.Worksheets.Item("registration_data")
$reg.Cells.Item(2, 1).Value = "Neeraja43@gmail.com"
$reg.Cells.Item(2, 2).Value = "Neeraja123@AP"

# --- 4. add filterorder_data at the end ----------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "filterorder_data"

$filterRows = @(
    @("email", "password", "username"),
    @("Rani12@gmail.com", "Neeru@34523456", "Rani12"),
    @("Rani13@gmail.com", "Neeru@34523456", "Rani13")
)
for ($r = 0; $r -lt $filterRows.Length; $r++) {
    $rowVals = $filterRows[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $newSheet.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}

Write-Host "edit complete"
